$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("Aceite") previously had D9 = 1 (numeric, percent-formatted).
# The author retyped it as literal text "100.0%" (leading apostrophe forces
# text entry while the cell keeps its existing percentage number format /
# picks up the quote-prefix flag), which also adds a new shared string.
$ws.Range("D9").Value = "'100.0%"

# H9 lost its explicit number-format styling entirely (back to the default
# "Normal" cell style), while its value (0) is unchanged.
$ws.Range("H9").Style = "Normal"

# The user ended up with D10 selected (e.g. pressed Enter after editing D9).
$ws.Range("D10").Select()
